# Add data for 2021-10-07
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-09-29"

# Update the September row label to reflect the new "through" date
$ws.Range("A10").Value = "September (through 09-29)"

# Update September row (row 10) values
$ws.Range("B10").Value = 30
$ws.Range("C10").Value = 44
$ws.Range("D10").Value = 74
$ws.Range("E10").Value = 54
$ws.Range("F10").Value = 71
$ws.Range("G10").Value = 112
$ws.Range("H10").Value = 172

# Update Total row (row 11) values
$ws.Range("B11").Value = 224
$ws.Range("C11").Value = 425
$ws.Range("D11").Value = 625
$ws.Range("E11").Value = 544
$ws.Range("F11").Value = 420
$ws.Range("G11").Value = 896
$ws.Range("H11").Value = 1242
